$d = $word.ActiveDocument

# The document contains a single inline picture (the generated QR code).
# It was resized (scaled up, aspect ratio preserved) from
# 2514600 x 2514600 EMU to 2933698 x 2933698 EMU.
#
# Word's InlineShape.Width / InlineShape.Height are expressed in points
# (1 pt = 12700 EMU), so convert the target EMU extents to points before
# assigning them. This updates both <wp:extent> and the picture's
# <a:ext> (pic:spPr/a:xfrm/a:ext) to the new size, matching the diff.
$emuPerPoint = 12700
$targetCx = 2933698
$targetCy = 2933698

$shape = $d.InlineShapes.Item(1)
$shape.Width = $targetCx / $emuPerPoint
$shape.Height = $targetCy / $emuPerPoint
